# [Add] SOS 포탄 구현 [Update] 풀링 적용 중
#
# Inserts a new "SOS" skill row into the SkillTable sheet (as row 13),
# pushing the existing GhostKnight (고스트나이트) and MiniFairy
# (미니페어리) summon rows down by one. Updates the FamiliarData sheet's
# SkillId references to the shifted skill ids, expands the SkillTable
# table range to include the new row, and updates the active
# sheet/selection state.

$wb = $excel.ActiveWorkbook

$wsSkill = $wb.Worksheets.Item("SkillTable")
$wsFam   = $wb.Worksheets.Item("FamiliarData")

# --- 1. Insert a new row above the old row 13 (고스트나이트), shifting
#        it (and 미니페어리 after it) down by one row. ---
$wsSkill.Rows.Item(13).Insert()

# --- 2. Expand the "표2" table on SkillTable so it covers the new row
#        (A3:M14 -> A3:M15). ---
$loSkill = $wsSkill.ListObjects.Item(1)
$loSkill.Resize($wsSkill.Range("A3:M15"))

# --- 3. Populate the newly-inserted row 13 with the SOS skill data. ---
$wsSkill.Range("A13").Value = 10010
$wsSkill.Range("B13").Value = "SOS"
$wsSkill.Range("C13").Value = "3개의 투사체가 일정 시간마다 날아갑니다."
$wsSkill.Range("D13").Value = "SPUM/0_Flat/Icon_Flat__18"
$wsSkill.Range("E13").Value = 1
$wsSkill.Range("F13").Value = 0.2
$wsSkill.Range("G13").Value = "14000, 14001, 14003, 14006"
$wsSkill.Range("H13").Value = 1.5
$wsSkill.Range("I13").Value = 18
$wsSkill.Range("J13").Value = 12
$wsSkill.Range("K13").Value = 1
$wsSkill.Range("L13").Value = "Skills/SOSBullet"
$wsSkill.Range("M13").Value = "Outside"

# --- 4. The 고스트나이트 / 미니페어리 rows (now rows 14 & 15) keep their
#        data, but their SkillId (column A) shifts up by one, since the
#        new SOS entry claimed id 10010. ---
$wsSkill.Range("A14").Value = 10011
$wsSkill.Range("A15").Value = 10012

# --- 5. FamiliarData references the SkillTable ids in column C; update
#        them to match the shifted ids above. ---
$wsFam.Range("C4").Value = 10011
$wsFam.Range("C5").Value = 10012

# --- 6. Update sheet selection / active tab state: SkillTable becomes
#        the active sheet (selection F24), FamiliarData is no longer the
#        active tab (selection C6). ---
$wsFam.Select() | Out-Null
$wsFam.Range("C6").Select() | Out-Null
$wsSkill.Select() | Out-Null
$wsSkill.Range("F24").Select() | Out-Null
